$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (rows 2-15), reordered and updated values, plus one new row (Greece)
$data = @(
    @("Germany", 52, 17.99),
    @("Sweden", 46, 15.92),
    @("United Kingdom", 43, 14.88),
    @("Ireland", 27, 9.34),
    @("France", 26, 9),
    @("Switzerland", 25, 8.65),
    @("Denmark", 16, 5.54),
    @("Netherlands", 13, 4.5),
    @("Spain", 12, 4.15),
    @("Italy", 9, 3.11),
    @("Norway", 8, 2.77),
    @("Belgium", 5, 1.73),
    @("Finland", 5, 1.73),
    @("Greece", 2, 0.6899999999999999)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# Apply the same style (bold font, border, centered/top alignment) as the other
# column-A cells to the newly added row's A cell by copying formats only.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
